$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "602.58") need to be
# forced to Text before assignment, otherwise Excel auto-converts the string
# to a numeric value (losing trailing zeros / introducing float rounding),
# which the source workbook never does (every Price cell is stored as text).
$textCells = @("D5", "D6", "D11", "D12", "D19", "D20", "D21", "D22", "D23", "D25", "D27", "D28", "D31", "D33", "D38", "D39", "D42", "D43", "D44", "D48", "D49")
foreach ($cellref in $textCells) {
    $ws.Range($cellref).NumberFormat = "@"
}

$ws.Range("D2").Value = '69.182.05'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '3.759.82'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '602.58'
$ws.Range("E5").Value = '  +0.22%  '
$ws.Range("D6").Value = '167.40'
$ws.Range("E6").Value = '  -0.58%  '
$ws.Range("D7").Value = '3.759.04'
$ws.Range("E7").Value = '  +0.36%  '
$ws.Range("E9").Value = '  +1.26%  '
$ws.Range("E10").Value = '  +3.49%  '
$ws.Range("D11").Value = '6.41'
$ws.Range("E11").Value = '  +2.07%  '
$ws.Range("D12").Value = '0.460'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("E14").Value = '  +1.92%  '
$ws.Range("D15").Value = '4.390.64'
$ws.Range("E15").Value = '  +0.59%  '
$ws.Range("D16").Value = '3.765.04'
$ws.Range("E16").Value = '  +0.81%  '
$ws.Range("D17").Value = '69.204.98'
$ws.Range("E17").Value = '  +0.38%  '
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("B19").Value = 'TRON'
$ws.Range("C19").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D19").Value = '0.113'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D20").Value = '17.26'
$ws.Range("E20").Value = '  +0.31%  '
$ws.Range("D21").Value = '11.22'
$ws.Range("E21").Value = '  +15.52%  '
$ws.Range("D22").Value = '494.08'
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").Value = '0.730'
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("E24").Value = '  +7.68%  '
$ws.Range("D25").Value = '84.95'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("D27").Value = '12.33'
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").Value = '10.11'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  +1.43%  '
$ws.Range("D31").Value = '8.18'
$ws.Range("E31").Value = '  +3.08%  '
$ws.Range("E32").Value = '  +2.31%  '
$ws.Range("D33").Value = '31.63'
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("D34").Value = '3.907.58'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").Value = '3.697.20'
$ws.Range("E35").Value = '  +0.68%  '
$ws.Range("E36").Value = '  -0.32%  '
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("B38").Value = 'Mantle'
$ws.Range("C38").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D38").Value = '1.02'
$ws.Range("E38").Value = '  +0.91%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '5.98'
$ws.Range("E39").Value = '  +3.42%  '
$ws.Range("E40").Value = '  +2.60%  '
$ws.Range("E41").Value = '  +0.96%  '
$ws.Range("D42").Value = '3.04'
$ws.Range("E42").Value = '  +4.87%  '
$ws.Range("D43").Value = '429.89'
$ws.Range("E43").Value = '  -1.78%  '
$ws.Range("D44").Value = '48.64'
$ws.Range("E44").Value = '  -0.76%  '
$ws.Range("E45").Value = '  +0.26%  '
$ws.Range("E46").Value = '  +0.99%  '
$ws.Range("E47").Value = '  -0.02%  '
$ws.Range("D48").Value = '40.47'
$ws.Range("E48").Value = '  -0.17%  '
$ws.Range("D49").Value = '141.20'
$ws.Range("E49").Value = '  -1.66%  '
$ws.Range("D50").Value = '2.797.86'
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("E51").Value = '  +0.54%  '

# Drop the temporary Text number-format again so these cells end up with no
# explicit style, matching the rest of the data cells in the sheet.
foreach ($cellref in $textCells) {
    $ws.Range($cellref).ClearFormats()
}
